$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4 (shifts old rows 4-8 down to 5-9),
# matching the "Azerbaijan Premier League" fixture being inserted ahead of "Danish 1st Division".
$ws.Rows.Item(4).Insert()

# Row 2 updates
$ws.Cells.Item(2,7).Value2 = 2.42
$ws.Cells.Item(2,14).Value2 = 4.1
$ws.Cells.Item(2,18).Value2 = 1.42
$ws.Cells.Item(2,24).Value2 = 980
$ws.Cells.Item(2,25).Value2 = 980
$ws.Cells.Item(2,26).Value2 = 980
$ws.Cells.Item(2,30).Value2 = 980
$ws.Cells.Item(2,31).Value2 = 980
$ws.Cells.Item(2,32).Value2 = 980
$ws.Cells.Item(2,34).Value2 = 980
$ws.Cells.Item(2,35).Value2 = 980
$ws.Cells.Item(2,36).Value2 = 980
$ws.Cells.Item(2,37).Value2 = 980
$ws.Cells.Item(2,38).Value2 = 980

# Row 3 updates
$ws.Cells.Item(3,9).Value2 = 2.48
$ws.Cells.Item(3,10).Value2 = 3.85
$ws.Cells.Item(3,11).Value2 = 3.95
$ws.Cells.Item(3,12).Value2 = 1.34
$ws.Cells.Item(3,16).Value2 = 2.24
$ws.Cells.Item(3,17).Value2 = 1.73
$ws.Cells.Item(3,18).Value2 = 1.48
$ws.Cells.Item(3,19).Value2 = 2.84
$ws.Cells.Item(3,22).Value2 = 1.67
$ws.Cells.Item(3,24).Value2 = 20
$ws.Cells.Item(3,26).Value2 = 20
$ws.Cells.Item(3,29).Value2 = 9
$ws.Cells.Item(3,32).Value2 = 24
$ws.Cells.Item(3,35).Value2 = 36
$ws.Cells.Item(3,38).Value2 = 1000
$ws.Cells.Item(3,41).Value2 = 16.5

# Row 4 updates
$ws.Cells.Item(4,1).Value2 = "Azerbaijan Premier League"
$ws.Cells.Item(4,2).NumberFormat = "@"
$ws.Cells.Item(4,2).Value2 = "2025-11-21"
$ws.Cells.Item(4,3).Value2 = "07:00:00"
$ws.Cells.Item(4,4).Value2 = "Karvan Evlakh"
$ws.Cells.Item(4,5).Value2 = "Kapaz Ganja"
$ws.Cells.Item(4,6).Value2 = 2.66
$ws.Cells.Item(4,7).Value2 = 3.75
$ws.Cells.Item(4,8).Value2 = 2.04
$ws.Cells.Item(4,9).Value2 = 2.9
$ws.Cells.Item(4,10).Value2 = 3
$ws.Cells.Item(4,11).Value2 = 7.4
$ws.Cells.Item(4,12).Value2 = 1.33
$ws.Cells.Item(4,13).Value2 = 1.01
$ws.Cells.Item(4,14).Value2 = 1.63
$ws.Cells.Item(4,15).Value2 = 1.01
$ws.Cells.Item(4,16).Value2 = 1.63
$ws.Cells.Item(4,17).Value2 = 1.85
$ws.Cells.Item(4,18).Value2 = 1.2
$ws.Cells.Item(4,19).Value2 = 1.85
$ws.Cells.Item(4,20).Value2 = 1.01
$ws.Cells.Item(4,21).Value2 = 1.01
$ws.Cells.Item(4,22).Value2 = 1.52
$ws.Cells.Item(4,23).Value2 = 1.36
$ws.Cells.Item(4,24).Value2 = 1000
$ws.Cells.Item(4,25).Value2 = 1000
$ws.Cells.Item(4,26).Value2 = 1000
$ws.Cells.Item(4,27).Value2 = 1000
$ws.Cells.Item(4,28).Value2 = 1000
$ws.Cells.Item(4,29).Value2 = 1000
$ws.Cells.Item(4,30).Value2 = 1000
$ws.Cells.Item(4,31).Value2 = 1000
$ws.Cells.Item(4,32).Value2 = 1000
$ws.Cells.Item(4,33).Value2 = 1000
$ws.Cells.Item(4,34).Value2 = 1000
$ws.Cells.Item(4,35).Value2 = 1000
$ws.Cells.Item(4,36).Value2 = 1000
$ws.Cells.Item(4,37).Value2 = 1000
$ws.Cells.Item(4,38).Value2 = 1000
$ws.Cells.Item(4,39).Value2 = 1000
$ws.Cells.Item(4,40).Value2 = 1000
$ws.Cells.Item(4,41).Value2 = 1000

# Row 5 updates
$ws.Cells.Item(5,6).Value2 = 1.68
$ws.Cells.Item(5,7).Value2 = 2.08
$ws.Cells.Item(5,8).Value2 = 1.92
$ws.Cells.Item(5,9).Value2 = 7.8
$ws.Cells.Item(5,10).Value2 = 3.65
$ws.Cells.Item(5,16).Value2 = 2.22

# Row 6 updates
$ws.Cells.Item(6,8).Value2 = 2.72
$ws.Cells.Item(6,10).Value2 = 3.8
$ws.Cells.Item(6,14).Value2 = 5.3
$ws.Cells.Item(6,17).Value2 = 1.66
$ws.Cells.Item(6,24).Value2 = 24
$ws.Cells.Item(6,27).Value2 = 42
$ws.Cells.Item(6,29).Value2 = 8.6
$ws.Cells.Item(6,31).Value2 = 27
$ws.Cells.Item(6,37).Value2 = 25
$ws.Cells.Item(6,39).Value2 = 1000
$ws.Cells.Item(6,40).Value2 = 16

# Row 7 updates
$ws.Cells.Item(7,11).Value2 = 4
$ws.Cells.Item(7,17).Value2 = 1.74
$ws.Cells.Item(7,27).Value2 = 23
$ws.Cells.Item(7,32).Value2 = 36
$ws.Cells.Item(7,33).Value2 = 19
$ws.Cells.Item(7,37).Value2 = 44

# Row 8 updates
$ws.Cells.Item(8,9).Value2 = 2
$ws.Cells.Item(8,10).Value2 = 3.4
$ws.Cells.Item(8,17).Value2 = 1.78

# Row 9 updates
$ws.Cells.Item(9,8).Value2 = 3.7
$ws.Cells.Item(9,17).Value2 = 2.02
$ws.Cells.Item(9,18).Value2 = 1.37
$ws.Cells.Item(9,34).Value2 = 18.5
$ws.Cells.Item(9,36).Value2 = 27
